# "sửa lịch sử một chút" - update match-history data on Sheet1.
#
# Adds four new header columns (C1:F1) describing the per-player match
# stats / history, fills in the new match-history values for the
# "anhlavodich" row (row 6: trận đấu/thắng/thua counts + history strings),
# and removes the stray leftover value in I8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1: new header cells C1:F1 --------------------------------------
# C1:E1 pick up the existing right-aligned "#,##0" column style (columns
# C-E already carry that style), F1 needs it applied explicitly since
# column F itself has no default style.
$ws.Range("C1").Value = "Số trận đấu"
$ws.Range("D1").Value = "Số lần thắng"
$ws.Range("E1").Value = "Số lần thua"
$ws.Range("F1").Value = "Lịch sử đấu (Thắng-1/Thua-0)"
$ws.Range("F1").NumberFormat = "#,##0"
$ws.Range("F1").HorizontalAlignment = -4152

# --- Row 6 ("anhlavodich"): updated stats + new history cells ----------
$ws.Range("B6").Value = 2200
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 2

$ws.Range("F6").Value = ";1;1"
$ws.Range("G6").Value = ";22;23"
$ws.Range("H6").Value = ";+100;+100"

# --- Row 8 ("taotaikhoan"): drop the stray leftover value in I8 --------
$ws.Range("I8").Clear()

# --- Column F is now a free-form history string, widen it so it's no
#     longer just an auto "best fit" column ---------------------------
$ws.Columns.Item(6).ColumnWidth = 26.3

# --- Leave the selection where the author left it on save --------------
$ws.Range("E15").Select() | Out-Null
